# "Modularized shopping cart class" - refresh the saved login/search test data
# on Sheet1: rotate the test account's email/password, turn the new email
# into a clickable mailto hyperlink, and leave Sheet1 as the active tab with
# the password cell selected (mirrors the author's last interactive state).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Swap out the old test credentials for new ones.
$ws1.Range("B2").Value = "audi.love25@gmail.com"
$ws1.Range("C2").Value = "Mitsubishi7!"

# Excel auto-hyperlinks e-mail addresses on entry; make it explicit here so
# the workbook gets a Hyperlink-styled cell + relationship for B2.
$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:audi.love25@gmail.com")

# Sheet1 becomes the active/visible sheet, selection parked on C2.
$ws1.Activate()
$ws1.Range("C2").Select()
